$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to be inserted into column B for each data row (2-24).
# Columns B..J shift right into C..K (old K value is dropped) to make
# room for this new "filter" column.
$newB = @{
    2  = 0.04231489763667018
    3  = 1.228061995268202
    4  = 0.5137754236260815
    5  = 0.4028038717171413
    6  = 0.5336388157440486
    7  = -1.102517691576566
    8  = 0.1737007515684039
    9  = 0.8389029408811082
    10 = -0.6954484448595206
    11 = -0.2465870357053012
    12 = -0.1975260465718366
    13 = 0.4425040297996861
    14 = -0.2720610750631522
    15 = -0.1065518669046048
    16 = -0.1895682054566924
    17 = 1.157000698704573
    18 = -0.4886691766355519
    19 = 1.10624937372658
    20 = -0.6446211617534254
    21 = -0.6387305113048862
    22 = 0.3668428211138005
    23 = -0.4578680368388337
    24 = 0.4181606776922825
}

for ($row = 2; $row -le 24; $row++) {
    # Shift existing values from column J (10) down to B (2) rightward one
    # column at a time, working right-to-left so we never clobber a value
    # before reading it. Old column K (11) value is discarded in the
    # process (shifted off the used range).
    for ($col = 10; $col -ge 2; $col--) {
        $srcVal = $ws.Cells.Item($row, $col).Value2
        $ws.Cells.Item($row, $col + 1).Value = $srcVal
    }

    # Write the new filter-column value into column B.
    $ws.Cells.Item($row, 2).Value = $newB[$row]
}
